# Updates cryptos list values (prices, volume %, and two coin-row swaps)
# per GitHub Actions scheduled data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.339.84"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "2.175.95"
$ws.Range("E3").Value = "  -2.06%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "252.73"
$ws.Range("E5").Value = "  +5.00%  "
$ws.Range("D6").Value = "0.606"
$ws.Range("E6").Value = "  -2.07%  "
$ws.Range("D7").Value = "73.23"
$ws.Range("E7").Value = "  -2.52%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.576"
$ws.Range("E9").Value = "  -4.57%  "
$ws.Range("D10").Value = "40.07"
$ws.Range("E10").Value = "  -2.93%  "
$ws.Range("D11").Value = "0.0912"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.101"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "6.73"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D14").Value = "2.504.13"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").Value = "14.12"
$ws.Range("E15").Value = "  -3.71%  "
$ws.Range("D16").Value = "2.156.24"
$ws.Range("E16").Value = "  -3.09%  "
$ws.Range("D17").Value = "0.767"
$ws.Range("E17").Value = "  -4.24%  "
$ws.Range("D18").Value = "42.251.44"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("E19").Value = "  -3.52%  "
$ws.Range("D20").Value = "70.55"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").Value = "5.85"
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("D22").Value = "226.03"
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "9.33"
$ws.Range("E23").Value = "  -7.42%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").Value = "2.10"
$ws.Range("E24").Value = "  -3.62%  "
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").Value = "10.42"
$ws.Range("E26").Value = "  -4.61%  "
$ws.Range("D27").Value = "3.38"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "2.16"
$ws.Range("E28").Value = "  -2.31%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.14"
$ws.Range("E29").Value = "  +2.26%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "170.80"
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "36.72"
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("D32").Value = "19.96"
$ws.Range("E32").Value = "  -1.77%  "
$ws.Range("D33").Value = "0.0809"
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("D34").Value = "5.08"
$ws.Range("E34").Value = "  -4.76%  "
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("D36").Value = "0.106"
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("D37").Value = "4.21"
$ws.Range("E37").Value = "  -4.24%  "
$ws.Range("D38").Value = "0.0334"
$ws.Range("E38").Value = "  +3.56%  "
$ws.Range("D39").Value = "11.70"
$ws.Range("E39").Value = "  -6.17%  "
$ws.Range("D40").Value = "2.04"
$ws.Range("E40").Value = "  -4.62%  "
$ws.Range("D41").Value = "0.195"
$ws.Range("E41").Value = "  -1.38%  "
$ws.Range("D42").Value = "58.89"
$ws.Range("E42").Value = "  -3.17%  "
$ws.Range("D43").Value = "5.10"
$ws.Range("E43").Value = "  -7.24%  "
$ws.Range("D44").Value = "101.31"
$ws.Range("E44").Value = "  +1.98%  "
$ws.Range("D45").Value = "2.46"
$ws.Range("E45").Value = "  +7.17%  "
$ws.Range("D46").Value = "0.0973"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("D47").Value = "0.461"
$ws.Range("E47").Value = "  +9.39%  "
$ws.Range("D48").Value = "8.19"
$ws.Range("E48").Value = "  -4.39%  "
$ws.Range("D49").Value = "1.08"
$ws.Range("E49").Value = "  -2.43%  "
$ws.Range("D50").Value = "1.12"
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("E51").Value = "  +0.08%  "

$wb.Save()
